$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph by scanning the
# document's paragraph collection (more reliable here than chaining through
# Range.Paragraphs on a Find result / collapsed range).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Docente\(s\) Responsável\(eis\)") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $para = $d.Paragraphs.Item($targetIndex)

    # Move to the end of that paragraph and insert a new paragraph right
    # after it with the instructor entry.
    $insertPoint = $para.Range
    $insertPoint.Collapse(0)
    $insertPoint.InsertAfter("5111420 - Talita Martins Lacerda`r")

    # The newly created paragraph is now the next one in the document;
    # give it the bullet-list style used for this kind of entry.
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Style = "ListBullet"
}
